$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold price/volume figures formatted as plain text (e.g. "17.00",
# "1.447.99", "0.01988"). Force the cell number format to Text ("@") before writing
# so Excel does not reinterpret/round them as numeric values.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.192.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.442.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9223"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -7.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.77"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.88%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3659"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.86%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3131"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.06"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.025"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06539"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.408"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.66"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.082"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.447.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.90%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9378"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05618"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.440"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.50"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.71%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.271"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.180.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.193"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "136.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.602.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.764"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8169"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.852"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07672"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.496"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06035"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.710"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.135"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.27"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.86%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01988"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.74%  "

# Row 41
$ws.Range("B41").Value = "Frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9398"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1829"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.106"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -15.07%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.519"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5162"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.771"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06344"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9933"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.88%  "
